# New submission synced: 2026-02-09 05:29:14
# Sheet "JSS 3C" - append a new Google-Form-style response row and
# normalize the previous row's "Admission No" value to a real number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3C")

# Row 4 ("Muhammad Ali zarami"): Admission No was stored as text "23";
# re-save it as the numeric value 23 (AI Score D4 stays 8, unchanged).
$ws.Cells.Item(4, 3).Value = 23

# Row 5: brand-new form submission.
$ws.Cells.Item(5, 1).Value = "2026-02-09 05:29:14"
$ws.Cells.Item(5, 2).Value = "Abubakar shettima mutawalli  "

# Admission No "1" is a single-digit numeric-looking value that the form
# sync kept as plain text (like C4 used to be). Force text storage, then
# drop the number-format override so no extra style sticks to the cell.
$ws.Cells.Item(5, 3).NumberFormat = "@"
$ws.Cells.Item(5, 3).Value = "1"
$ws.Cells.Item(5, 3).ClearFormats()

$ws.Cells.Item(5, 4).Value = 8
